$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (changed) date column C for every existing
#    data row (2..471) from 45189 to 45190.
for ($r = 2; $r -le 471; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}

# 2) Row 471 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(471).RowHeight = 15

# 3) Append four new data rows (472-475).
$newRows = @(
    @{ Row = 472; A = "A 44249-2023"; B = 45188; C = 45190; G = 0.8 },
    @{ Row = 473; A = "A 44211-2023"; B = 45188; C = 45190; G = 1.2 },
    @{ Row = 474; A = "A 44243-2023"; B = 45188; C = 45190; G = 0.8 },
    @{ Row = 475; A = "A 44435-2023"; B = 45189; C = 45190; G = 0.8 }
)

foreach ($row in $newRows) {
    $r = $row.Row

    $ws.Cells.Item($r, 1).Value = $row.A

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = "VÄSTRA GÖTALANDS LÄN"
    $ws.Cells.Item($r, 5).Value = "MARK"

    $ws.Cells.Item($r, 7).Value = $row.G

    for ($c = 8; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }

    # R column: empty cell but carrying the wrap-text style (s="2"),
    # matching the existing rows above it.
    $ws.Cells.Item($r, 18).WrapText = $true

    # Rows 472-474 pick up an explicit row height; row 475 does not.
    if ($r -le 474) {
        $ws.Rows.Item($r).RowHeight = 15
    }
}
